# Update column G ("K" = strikeouts) values on the active worksheet.
# These replace the previous "Strike#" counts with actual strikeout (K) totals,
# per the commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 10
    3  = 14
    4  = 5
    5  = 6
    6  = 8
    7  = 10
    8  = 11
    9  = 8
    10 = 9
    11 = 6
    12 = 9
    13 = 15
    14 = 14
    15 = 14
    16 = 7
    17 = 5
    18 = 6
    19 = 7
    20 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
